# Auto-generated script to update cryptos price/volume columns (D, E) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.926.86"
$ws.Range("E2").Value = "  -4.67%  "
$ws.Range("D3").Value = "1.738.02"
$ws.Range("E3").Value = "  -5.04%  "
$ws.Range("E4").Value = "  -0.24%  "
$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "226.37"
$c.Style = $origStyle
$ws.Range("E5").Value = "  -3.89%  "
$c = $ws.Range("D6")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.5788"
$c.Style = $origStyle
$ws.Range("E6").Value = "  -3.99%  "
$c = $ws.Range("D7")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = $origStyle
$ws.Range("E7").Value = "  -0.24%  "
$c = $ws.Range("D8")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.2726"
$c.Style = $origStyle
$ws.Range("E8").Value = "  -1.64%  "
$c = $ws.Range("D9")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "23.26"
$c.Style = $origStyle
$ws.Range("E9").Value = "  -1.41%  "
$c = $ws.Range("D10")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.06602"
$c.Style = $origStyle
$ws.Range("E10").Value = "  -5.36%  "
$c = $ws.Range("D11")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.07557"
$c.Style = $origStyle
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "1.741.37"
$ws.Range("E12").Value = "  -5.00%  "
$c = $ws.Range("D13")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.701"
$c.Style = $origStyle
$ws.Range("E13").Value = "  -1.15%  "
$c = $ws.Range("D14")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.6017"
$c.Style = $origStyle
$ws.Range("E14").Value = "  -4.89%  "
$ws.Range("D15").Value = "1.974.67"
$ws.Range("E15").Value = "  -5.05%  "
$c = $ws.Range("D16")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "74.55"
$c.Style = $origStyle
$ws.Range("E16").Value = "  -4.26%  "
$c = $ws.Range("D17")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.000008719"
$c.Style = $origStyle
$ws.Range("E17").Value = "  -11.57%  "
$ws.Range("D18").Value = "27.934.32"
$ws.Range("E18").Value = "  -3.64%  "
$c = $ws.Range("D19")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.315"
$c.Style = $origStyle
$ws.Range("E19").Value = "  -5.00%  "
$ws.Range("E20").Value = "  -0.21%  "
$c = $ws.Range("D21")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "205.41"
$c.Style = $origStyle
$ws.Range("E21").Value = "  -5.74%  "
$c = $ws.Range("D22")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "11.27"
$c.Style = $origStyle
$ws.Range("E22").Value = "  -2.81%  "
$c = $ws.Range("D23")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.617"
$c.Style = $origStyle
$ws.Range("E23").Value = "  -4.28%  "
$ws.Range("E24").Value = "  -0.23%  "
$c = $ws.Range("D25")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "150.06"
$c.Style = $origStyle
$ws.Range("E25").Value = "  -4.19%  "
$c = $ws.Range("D26")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.132"
$c.Style = $origStyle
$ws.Range("E26").Value = "  +1.84%  "
$c = $ws.Range("D27")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.1230"
$c.Style = $origStyle
$ws.Range("E27").Value = "  -4.90%  "
$c = $ws.Range("D28")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "16.13"
$c.Style = $origStyle
$ws.Range("E28").Value = "  -2.60%  "
$c = $ws.Range("D29")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.384"
$c.Style = $origStyle
$ws.Range("E29").Value = "  -2.73%  "
$c = $ws.Range("D30")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.06150"
$c.Style = $origStyle
$ws.Range("E30").Value = "  -4.41%  "
$c = $ws.Range("D31")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.390"
$c.Style = $origStyle
$c = $ws.Range("D32")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.734"
$c.Style = $origStyle
$ws.Range("E32").Value = "  -2.74%  "
$c = $ws.Range("D33")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.723"
$c.Style = $origStyle
$ws.Range("E33").Value = "  -2.08%  "
$ws.Range("E34").Value = "  -3.63%  "
$ws.Range("E35").Value = "  -5.70%  "
$c = $ws.Range("D36")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.6414"
$c.Style = $origStyle
$ws.Range("E36").Value = "  -1.13%  "
$c = $ws.Range("D37")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.417"
$c.Style = $origStyle
$ws.Range("E37").Value = "  -5.05%  "
$c = $ws.Range("D38")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.717"
$c.Style = $origStyle
$ws.Range("E38").Value = "  -1.47%  "
$c = $ws.Range("D39")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.01668"
$c.Style = $origStyle
$ws.Range("E39").Value = "  -5.09%  "
$ws.Range("D40").Value = "1.128.50"
$ws.Range("E40").Value = "  -1.33%  "
$c = $ws.Range("D41")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.172"
$c.Style = $origStyle
$ws.Range("E41").Value = "  -6.65%  "
$c = $ws.Range("D42")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.8741"
$c.Style = $origStyle
$ws.Range("E42").Value = "  -2.24%  "
$c = $ws.Range("D43")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = $origStyle
$ws.Range("E43").Value = "  -0.09%  "
$c = $ws.Range("D44")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "99.65"
$c.Style = $origStyle
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("D45").Value = "1.888.66"
$ws.Range("E45").Value = "  -5.33%  "
$c = $ws.Range("D46")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "59.31"
$c.Style = $origStyle
$ws.Range("E46").Value = "  -4.81%  "
$c = $ws.Range("D47")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.577"
$c.Style = $origStyle
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("E48").Value = "  -6.04%  "
$c = $ws.Range("D49")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.248"
$c.Style = $origStyle
$ws.Range("E49").Value = "  -3.18%  "
$c = $ws.Range("D50")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.05378"
$c.Style = $origStyle
$ws.Range("E50").Value = "  -2.20%  "
$c = $ws.Range("D51")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.4415"
$c.Style = $origStyle
$ws.Range("E51").Value = "  -2.81%  "
